# Update the cryptocurrency price/hour table on the active worksheet.
# - Column D ("Price") gets new values for the rows listed below.
# - Column G ("Hora") changes from "10" to "11" for every data row (2-51).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Price (column D) values, keyed by row number.
$priceUpdates = @{
    2  = "275.25"
    3  = "22.97"
    4  = "6.322"
    5  = "0.06233"
    6  = "3.656"
    7  = "6.707"
    8  = "1.369"
    9  = "0.8332"
    12 = "0.08316"
    13 = "0.03348"
    14 = "0.03101"
    15 = "0.09326"
    16 = "3.889"
    17 = "0.001647"
    18 = "0.04782"
    19 = "0.006223"
    20 = "0.005569"
    21 = "0.001088"
    23 = "3.729"
    24 = "2.382"
    25 = "0.3385"
    26 = "0.1270"
    40 = "0.04700"
    41 = "0.007027"
    42 = "0.1168"
    44 = "0.01221"
    45 = "0.00006267"
    48 = "0.03061"
}

foreach ($row in $priceUpdates.Keys) {
    $cell = $ws.Range("D$row")
    # Force text storage (matches the workbook's existing inline-string
    # cells) rather than letting Excel auto-detect these as numbers.
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$row]
    $cell.Style = "Normal"
}

# Column G ("Hora") moves from 10 to 11 for every data row (2 through 51).
for ($row = 2; $row -le 51; $row++) {
    $cell = $ws.Range("G$row")
    $cell.NumberFormat = "@"
    $cell.Value = "11"
    $cell.Style = "Normal"
}
